$wb = $excel.ActiveWorkbook

# --- Sheet "chains": add a new chain row for "heat" ---
$chains = $wb.Worksheets.Item("chains")
$chains.Range("A9").Value = "heat"
$chains.Range("B9").Value = "heat"
$chains.Range("C9").Value = "outflow"
$chains.Range("E9").Value = "simple_heat"
$chains.Range("E9").NumberFormat = "@"

# --- Sheet "connections": connect the heat use (waste heat from power -> CO2
#     Capture's recovered heat; CO2 Capture's heat -> new "heat" chain) ---
$conn = $wb.Worksheets.Item("connections")

# Make room for two new connection rows, right after the existing
# "CO2 Capture inflow electricity" row. Excel shifts the "CO2 Storage inflow
# electricity" row down to row 17 automatically (inheriting formatting).
$conn.Rows("15:16").Insert()

# Row 13: CO2 Capture <- electricity (was row 14; now reordered to row 13)
$conn.Range("B13").Value = "CO2 Capture"
$conn.Range("C13").Value = "simple_CO2capture"
$conn.Range("D13").Value = "inflow"
$conn.Range("E13").Value = "electricity"
$conn.Range("F13").Value = "electricity"
$conn.Range("G13").Value = "outflow"
$conn.Range("H13").Value = "simple_power"
$conn.Range("I13").Value = "power"

# Row 14 (new): power's waste heat -> CO2 Capture's recovered heat, replacing heat
$conn.Range("B14").Value = "power"
$conn.Range("C14").Value = "simple_power"
$conn.Range("C14").NumberFormat = "@"
$conn.Range("D14").Value = "outflow"
$conn.Range("E14").Value = "waste heat"
$conn.Range("F14").Value = "recovered heat"
$conn.Range("G14").Value = "inflows"
$conn.Range("H14").Value = "simple_CO2capture"
$conn.Range("I14").Value = "CO2 Capture"
$conn.Range("J14").Value = "heat"

# Row 15 (new): CO2 Capture's heat inflow -> new "heat" chain outflow
$conn.Range("B15").Value = "CO2 Capture"
$conn.Range("C15").Value = "simple_CO2capture"
$conn.Range("D15").Value = "inflow"
$conn.Range("E15").Value = "heat"
$conn.Range("F15").Value = "heat"
$conn.Range("G15").Value = "outflows"
$conn.Range("H15").Value = "simple_heat"
$conn.Range("I15").Value = "heat"

# Row 16: CO2 Capture -> CO2 Storage (was row 13; now reordered to row 16)
$conn.Range("B16").Value = "CO2 Capture"
$conn.Range("C16").Value = "simple_CO2capture"
$conn.Range("D16").Value = "outflow"
$conn.Range("E16").Value = "compressed CO2"
$conn.Range("F16").Value = "compressed CO2"
$conn.Range("G16").Value = "inflows"
$conn.Range("H16").Value = "simple_CO2storage"
$conn.Range("I16").Value = "CO2 Storage"

# Row 17 (was row 15 "CO2 Storage inflow electricity") is already correct
# after the row insert/shift above.

# Leave the selection on "connections" where the new rows were added, then
# finish back on "chains" (which becomes the active tab again).
[void]$conn.Range("I15").Select()
[void]$chains.Activate()
[void]$chains.Range("A13").Select()
